# DAS-677 - CCRU - Creation of Scenes for SOVI SOCVI
#
# The "Canteen" sheet lists KPIs together with the Scenes (columns X =
# "Scenes to include", Y = "Scenes to exclude") that each KPI applies to.
# This edit adds the newly created SOVI/SOCVI scene names to the existing
# Scene lists used by several KPI rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Scenes to exclude" (column Y): every row that excluded the "Panoramic
# Photo" scene now also excludes the new "SS_Panoramic Photo" scene.
$panoramicPhotoRows = @(4,5,6,7,8,9,10,11,12,13,14,16,17,18,20,21,22,23,25,26,27,28,29,30,31)
foreach ($r in $panoramicPhotoRows) {
    $ws.Range("Y$r").Value = "Panoramic Photo, SS_Panoramic Photo"
}

# "Scenes to include" (column X): the Cooler Prime Position KPI (row 38)
# now also includes the new "SS_Panoramic photo of Cooler - Horeca" scene.
$ws.Range("X38").Value = "Panoramic photo of Cooler, SS_Panoramic photo of Cooler - Horeca"

# "Scenes to include" (column X): the Activation KPIs (rows 43-49) now
# reference the renamed/expanded Menu Board & Cash Zone scenes.
$menuBoardRows = @(43,44,45,46,47,48,49)
foreach ($r in $menuBoardRows) {
    $ws.Range("X$r").Value = "Menu Board, Cash Zone, SS_Cash Zone - Canteen, QSR, SS_Menu Board - Canteen, QSR"
}

# Reflect the author's final cell selection on the sheet.
$ws.Range("AD5").Select()
